$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1932773109243698
$ws.Range("C2").Value = 0.5294117647058824
$ws.Range("P2").Value = 0.1974789915966386
$ws.Range("S2").Value = 0.07983193277310924
$ws.Range("C3").Value = 0.03007518796992481
$ws.Range("J3").Value = 0.01503759398496241
$ws.Range("P3").Value = 0.7744360902255639
$ws.Range("S3").Value = 0.1804511278195489
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2647058823529412
$ws.Range("B6").Value = 0.04838709677419355
$ws.Range("F6").Value = 0.03763440860215054
$ws.Range("J6").Value = 0.2849462365591398
$ws.Range("O6").Value = 0.02150537634408602
$ws.Range("Q6").Value = 0.1505376344086022
$ws.Range("R6").Value = 0.05376344086021505
$ws.Range("S6").Value = 0.4032258064516129
$ws.Range("B7").Value = 0.1
$ws.Range("D7").Value = 0.00909090909090909
$ws.Range("F7").Value = 0.04090909090909091
$ws.Range("J7").Value = 0.1318181818181818
$ws.Range("O7").Value = 0.01363636363636364
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.05454545454545454
$ws.Range("S7").Value = 0.45
$ws.Range("B8").Value = 0.08314606741573034
$ws.Range("D8").Value = 0.008988764044943821
$ws.Range("F8").Value = 0.04044943820224719
$ws.Range("J8").Value = 0.1146067415730337
$ws.Range("O8").Value = 0.02247191011235955
$ws.Range("Q8").Value = 0.1865168539325843
$ws.Range("R8").Value = 0.06067415730337079
$ws.Range("S8").Value = 0.4831460674157304
$ws.Range("B9").Value = 0.0728476821192053
$ws.Range("F9").Value = 0.05298013245033113
$ws.Range("J9").Value = 0.1258278145695364
$ws.Range("O9").Value = 0.01324503311258278
$ws.Range("Q9").Value = 0.1854304635761589
$ws.Range("R9").Value = 0.06622516556291391
$ws.Range("S9").Value = 0.4834437086092715
$ws.Range("B10").Value = 0.09927797833935018
$ws.Range("D10").Value = 0.02617328519855596
$ws.Range("F10").Value = 0.07310469314079422
$ws.Range("J10").Value = 0.1254512635379061
$ws.Range("O10").Value = 0.02075812274368231
$ws.Range("Q10").Value = 0.2536101083032491
$ws.Range("R10").Value = 0.05144404332129964
$ws.Range("S10").Value = 0.3501805054151624
$ws.Range("G11").Value = 0.1783625730994152
$ws.Range("J11").Value = 0.09649122807017543
$ws.Range("K11").Value = 0.2251461988304093
$ws.Range("L11").Value = 0.4853801169590643
$ws.Range("S11").Value = 0.01461988304093567
$ws.Range("G12").Value = 0.7738095238095238
$ws.Range("J12").Value = 0.1547619047619048
$ws.Range("S12").Value = 0.07142857142857142
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.0160427807486631
$ws.Range("H15").Value = 0.1122994652406417
$ws.Range("I15").Value = 0.0374331550802139
$ws.Range("J15").Value = 0.3689839572192513
$ws.Range("K15").Value = 0.06417112299465241
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.3422459893048128
$ws.Range("F16").Value = 0.03488372093023256
$ws.Range("H16").Value = 0.1104651162790698
$ws.Range("I16").Value = 0.05232558139534884
$ws.Range("J16").Value = 0.4186046511627907
$ws.Range("K16").Value = 0.1395348837209302
$ws.Range("M16").Value = 0.02325581395348837
$ws.Range("O16").Value = 0.02906976744186046
$ws.Range("S16").Value = 0.1918604651162791
$ws.Range("F17").Value = 0.01754385964912281
$ws.Range("H17").Value = 0.1842105263157895
$ws.Range("I17").Value = 0.09210526315789473
$ws.Range("J17").Value = 0.375
$ws.Range("K17").Value = 0.1206140350877193
$ws.Range("M17").Value = 0.02192982456140351
$ws.Range("O17").Value = 0.04605263157894737
$ws.Range("S17").Value = 0.1425438596491228
$ws.Range("F18").Value = 0.04385964912280702
$ws.Range("H18").Value = 0.1754385964912281
$ws.Range("I18").Value = 0.03508771929824561
$ws.Range("J18").Value = 0.412280701754386
$ws.Range("K18").Value = 0.131578947368421
$ws.Range("M18").Value = 0.03508771929824561
$ws.Range("O18").Value = 0.07894736842105263
$ws.Range("S18").Value = 0.08771929824561403
$ws.Range("F19").Value = 0.01837060702875399
$ws.Range("H19").Value = 0.2476038338658147
$ws.Range("I19").Value = 0.07268370607028754
$ws.Range("J19").Value = 0.3194888178913738
$ws.Range("K19").Value = 0.1222044728434505
$ws.Range("M19").Value = 0.01996805111821086
$ws.Range("N19").Value = 0.0007987220447284345
$ws.Range("O19").Value = 0.06070287539936102
$ws.Range("S19").Value = 0.1381789137380192
